# Updates cached market-price / profit figures (columns H-N) on several
# rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets,
# as refreshed by the scheduled data-update runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4749  # ALC!H64
$ws.Cells.Item(64, 9).Value = 3700  # ALC!I64
$ws.Cells.Item(64, 10).Value = 5198.5713  # ALC!J64
$ws.Cells.Item(64, 11).Value = 3700  # ALC!K64
$ws.Cells.Item(64, 12).Value = 5198.5713  # ALC!L64
$ws.Cells.Item(64, 13).Value = -3452  # ALC!M64
$ws.Cells.Item(64, 14).Value = -5694.5713  # ALC!N64

$ws.Cells.Item(67, 8).Value = 4749  # ALC!H67
$ws.Cells.Item(67, 9).Value = 3700  # ALC!I67
$ws.Cells.Item(67, 10).Value = 5198.5713  # ALC!J67
$ws.Cells.Item(67, 11).Value = 3700  # ALC!K67
$ws.Cells.Item(67, 12).Value = 5198.5713  # ALC!L67
$ws.Cells.Item(67, 13).Value = -2842  # ALC!M67
$ws.Cells.Item(67, 14).Value = -6914.5713  # ALC!N67

$ws.Cells.Item(113, 8).Value = 3102.8572  # ALC!H113
$ws.Cells.Item(113, 9).Value = 3084  # ALC!I113
$ws.Cells.Item(113, 10).Value = 3150  # ALC!J113
$ws.Cells.Item(113, 11).Value = 3084  # ALC!K113
$ws.Cells.Item(113, 12).Value = 3150  # ALC!L113
$ws.Cells.Item(113, 13).Value = 170  # ALC!M113
$ws.Cells.Item(113, 14).Value = -9658  # ALC!N113

$ws.Cells.Item(121, 8).Value = 1273.5  # ALC!H121
$ws.Cells.Item(121, 10).Value = 1273.5  # ALC!J121
$ws.Cells.Item(121, 12).Value = 3820.5  # ALC!L121
$ws.Cells.Item(121, 14).Value = -7314.5  # ALC!N121

$ws.Cells.Item(129, 8).Value = 973.4  # ALC!H129
$ws.Cells.Item(129, 9).Value = 602.5  # ALC!I129
$ws.Cells.Item(129, 10).Value = 1021.25806  # ALC!J129
$ws.Cells.Item(129, 11).Value = 1807.5  # ALC!K129
$ws.Cells.Item(129, 12).Value = 3063.77418  # ALC!L129
$ws.Cells.Item(129, 13).Value = 3192.5  # ALC!M129
$ws.Cells.Item(129, 14).Value = -13063.77418  # ALC!N129

$ws.Cells.Item(138, 8).Value = 2224571.2  # ALC!H138
$ws.Cells.Item(138, 9).Value = 691.5  # ALC!I138
$ws.Cells.Item(138, 10).Value = 3879551.5  # ALC!J138
$ws.Cells.Item(138, 11).Value = 2074.5  # ALC!K138
$ws.Cells.Item(138, 12).Value = 11638654.5  # ALC!L138
$ws.Cells.Item(138, 13).Value = 3065.5  # ALC!M138
$ws.Cells.Item(138, 14).Value = -11648934.5  # ALC!N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 143144270  # ARM!H61
$ws.Cells.Item(61, 9).Value = 200200980  # ARM!I61
$ws.Cells.Item(61, 10).Value = 502500  # ARM!J61
$ws.Cells.Item(61, 11).Value = 200200980  # ARM!K61
$ws.Cells.Item(61, 12).Value = 502500  # ARM!L61
$ws.Cells.Item(61, 13).Value = -200200768  # ARM!M61
$ws.Cells.Item(61, 14).Value = -502924  # ARM!N61

$ws.Cells.Item(136, 8).Value = 143144270  # ARM!H136
$ws.Cells.Item(136, 9).Value = 200200980  # ARM!I136
$ws.Cells.Item(136, 10).Value = 502500  # ARM!J136
$ws.Cells.Item(136, 11).Value = 600602940  # ARM!K136
$ws.Cells.Item(136, 12).Value = 1507500  # ARM!L136
$ws.Cells.Item(136, 13).Value = -600600390  # ARM!M136
$ws.Cells.Item(136, 14).Value = -1512600  # ARM!N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4786.207  # BSM!H134
$ws.Cells.Item(134, 9).Value = 4714.2856  # BSM!I134
$ws.Cells.Item(134, 10).Value = 4975  # BSM!J134
$ws.Cells.Item(134, 11).Value = 14142.8568  # BSM!K134
$ws.Cells.Item(134, 12).Value = 14925  # BSM!L134
$ws.Cells.Item(134, 13).Value = -11607.8568  # BSM!M134
$ws.Cells.Item(134, 14).Value = -19995  # BSM!N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2282.1606  # CRP!H31
$ws.Cells.Item(31, 9).Value = 1405.561  # CRP!I31
$ws.Cells.Item(31, 10).Value = 4678.2  # CRP!J31
$ws.Cells.Item(31, 11).Value = 1405.561  # CRP!K31
$ws.Cells.Item(31, 12).Value = 4678.2  # CRP!L31
$ws.Cells.Item(31, 13).Value = -1110.561  # CRP!M31
$ws.Cells.Item(31, 14).Value = -5268.2  # CRP!N31

$ws.Cells.Item(34, 8).Value = 2282.1606  # CRP!H34
$ws.Cells.Item(34, 9).Value = 1405.561  # CRP!I34
$ws.Cells.Item(34, 10).Value = 4678.2  # CRP!J34
$ws.Cells.Item(34, 11).Value = 1405.561  # CRP!K34
$ws.Cells.Item(34, 12).Value = 4678.2  # CRP!L34
$ws.Cells.Item(34, 13).Value = -1203.561  # CRP!M34
$ws.Cells.Item(34, 14).Value = -5082.2  # CRP!N34

$ws.Cells.Item(58, 8).Value = 55558004  # CRP!H58
$ws.Cells.Item(58, 9).Value = 76924850  # CRP!I58
$ws.Cells.Item(58, 10).Value = 4200.2  # CRP!J58
$ws.Cells.Item(58, 11).Value = 76924850  # CRP!K58
$ws.Cells.Item(58, 12).Value = 4200.2  # CRP!L58
$ws.Cells.Item(58, 13).Value = -76924647  # CRP!M58
$ws.Cells.Item(58, 14).Value = -4606.2  # CRP!N58

$ws.Cells.Item(136, 8).Value = 55558004  # CRP!H136
$ws.Cells.Item(136, 9).Value = 76924850  # CRP!I136
$ws.Cells.Item(136, 10).Value = 4200.2  # CRP!J136
$ws.Cells.Item(136, 11).Value = 230774550  # CRP!K136
$ws.Cells.Item(136, 12).Value = 12600.6  # CRP!L136
$ws.Cells.Item(136, 13).Value = -230772000  # CRP!M136
$ws.Cells.Item(136, 14).Value = -17700.6  # CRP!N136

$ws.Cells.Item(138, 8).Value = 42000  # CRP!H138
$ws.Cells.Item(138, 10).Value = 42000  # CRP!J138
$ws.Cells.Item(138, 12).Value = 42000  # CRP!L138
$ws.Cells.Item(138, 14).Value = -52280  # CRP!N138

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 979.4  # CUL!H68
$ws.Cells.Item(68, 9).Value = 798  # CUL!I68
$ws.Cells.Item(68, 10).Value = 1251.5  # CUL!J68
$ws.Cells.Item(68, 11).Value = 2394  # CUL!K68
$ws.Cells.Item(68, 12).Value = 3754.5  # CUL!L68
$ws.Cells.Item(68, 13).Value = -1583  # CUL!M68
$ws.Cells.Item(68, 14).Value = -5376.5  # CUL!N68

$ws.Cells.Item(71, 8).Value = 979.4  # CUL!H71
$ws.Cells.Item(71, 9).Value = 798  # CUL!I71
$ws.Cells.Item(71, 10).Value = 1251.5  # CUL!J71
$ws.Cells.Item(71, 11).Value = 7182  # CUL!K71
$ws.Cells.Item(71, 12).Value = 11263.5  # CUL!L71
$ws.Cells.Item(71, 13).Value = -3126  # CUL!M71
$ws.Cells.Item(71, 14).Value = -19375.5  # CUL!N71

$ws.Cells.Item(94, 8).Value = 3270.7693  # CUL!H94
$ws.Cells.Item(94, 10).Value = 3809.0908  # CUL!J94
$ws.Cells.Item(94, 12).Value = 11427.2724  # CUL!L94
$ws.Cells.Item(94, 14).Value = -12779.2724  # CUL!N94

$ws.Cells.Item(131, 8).Value = 1340.5883  # CUL!H131
$ws.Cells.Item(131, 10).Value = 1396.129  # CUL!J131
$ws.Cells.Item(131, 12).Value = 4188.387  # CUL!L131
$ws.Cells.Item(131, 14).Value = -14268.387  # CUL!N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 20908.666  # GSM!H136
$ws.Cells.Item(136, 10).Value = 20908.666  # GSM!J136
$ws.Cells.Item(136, 12).Value = 62725.99800000001  # GSM!L136
$ws.Cells.Item(136, 14).Value = -67825.998  # GSM!N136

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2699.7778  # LTW!H7
$ws.Cells.Item(7, 9).Value = 2233.3333  # LTW!I7
$ws.Cells.Item(7, 10).Value = 2933  # LTW!J7
$ws.Cells.Item(7, 11).Value = 2233.3333  # LTW!K7
$ws.Cells.Item(7, 12).Value = 2933  # LTW!L7
$ws.Cells.Item(7, 13).Value = -2121.3333  # LTW!M7
$ws.Cells.Item(7, 14).Value = -3157  # LTW!N7

$ws.Cells.Item(22, 8).Value = 700.1429  # LTW!H22
$ws.Cells.Item(22, 9).Value = 700.1429  # LTW!I22
$ws.Cells.Item(22, 10).Value = 0  # LTW!J22
$ws.Cells.Item(22, 11).Value = 700.1429  # LTW!K22
$ws.Cells.Item(22, 12).Value = 0  # LTW!L22
$ws.Cells.Item(22, 13).Value = -405.1429000000001  # LTW!M22
$ws.Cells.Item(22, 14).ClearContents()  # LTW!N22

$ws.Cells.Item(27, 8).Value = 700.1429  # LTW!H27
$ws.Cells.Item(27, 9).Value = 700.1429  # LTW!I27
$ws.Cells.Item(27, 10).Value = 0  # LTW!J27
$ws.Cells.Item(27, 11).Value = 700.1429  # LTW!K27
$ws.Cells.Item(27, 12).Value = 0  # LTW!L27
$ws.Cells.Item(27, 13).Value = -593.1429  # LTW!M27
$ws.Cells.Item(27, 14).ClearContents()  # LTW!N27

$ws.Cells.Item(40, 8).Value = 2489.5  # LTW!H40
$ws.Cells.Item(40, 9).Value = 2477.2222  # LTW!I40
$ws.Cells.Item(40, 10).Value = 2600  # LTW!J40
$ws.Cells.Item(40, 11).Value = 2477.2222  # LTW!K40
$ws.Cells.Item(40, 12).Value = 2600  # LTW!L40
$ws.Cells.Item(40, 13).Value = -2341.2222  # LTW!M40
$ws.Cells.Item(40, 14).Value = -2872  # LTW!N40

$ws.Cells.Item(122, 8).Value = 2883.45  # LTW!H122
$ws.Cells.Item(122, 9).Value = 2256  # LTW!I122
$ws.Cells.Item(122, 10).Value = 3396.818  # LTW!J122
$ws.Cells.Item(122, 11).Value = 6768  # LTW!K122
$ws.Cells.Item(122, 12).Value = 10190.454  # LTW!L122
$ws.Cells.Item(122, 13).Value = -4318  # LTW!M122
$ws.Cells.Item(122, 14).Value = -15090.454  # LTW!N122

$ws.Cells.Item(126, 8).Value = 2699.7778  # LTW!H126
$ws.Cells.Item(126, 9).Value = 2233.3333  # LTW!I126
$ws.Cells.Item(126, 10).Value = 2933  # LTW!J126
$ws.Cells.Item(126, 11).Value = 6699.999899999999  # LTW!K126
$ws.Cells.Item(126, 12).Value = 8799  # LTW!L126
$ws.Cells.Item(126, 13).Value = -4229.999899999999  # LTW!M126
$ws.Cells.Item(126, 14).Value = -13739  # LTW!N126

$ws.Cells.Item(132, 8).Value = 38781.707  # LTW!H132
$ws.Cells.Item(132, 9).Value = 17382.719  # LTW!I132
$ws.Cells.Item(132, 10).Value = 114867  # LTW!J132
$ws.Cells.Item(132, 11).Value = 52148.15700000001  # LTW!K132
$ws.Cells.Item(132, 12).Value = 344601  # LTW!L132
$ws.Cells.Item(132, 13).Value = -49618.15700000001  # LTW!M132
$ws.Cells.Item(132, 14).Value = -349661  # LTW!N132

$ws.Cells.Item(136, 8).Value = 37006.156  # LTW!H136
$ws.Cells.Item(136, 9).Value = 21838.771  # LTW!I136
$ws.Cells.Item(136, 11).Value = 65516.313  # LTW!K136
$ws.Cells.Item(136, 13).Value = -62966.313  # LTW!M136

$ws.Cells.Item(140, 8).Value = 55264.5  # LTW!H140
$ws.Cells.Item(140, 10).Value = 55264.5  # LTW!J140
$ws.Cells.Item(140, 12).Value = 55264.5  # LTW!L140
$ws.Cells.Item(140, 14).Value = -65624.5  # LTW!N140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 439.14285  # WVR!H107
$ws.Cells.Item(107, 9).Value = 445.66666  # WVR!I107
$ws.Cells.Item(107, 11).Value = 1336.99998  # WVR!K107
$ws.Cells.Item(107, 13).Value = 583.00002  # WVR!M107

